$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 5 (pushes the existing row 5 and everything
# below it down by one), for the new weekly price observation dated 44677.
$ws.Rows("5:5").Insert()
$ws.Cells.Item(5, 1).Value = 10
$ws.Cells.Item(5, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(5, 3).Value = "La Araucanía"
$ws.Cells.Item(5, 4).Value = 44677
$ws.Cells.Item(5, 5).Value = 9
$ws.Cells.Item(5, 6).Value = 100112042
$ws.Cells.Item(5, 7).Value = "Locoto"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 20
$ws.Cells.Item(5, 11).Value = 5500
$ws.Cells.Item(5, 12).Value = 5500
$ws.Cells.Item(5, 13).Value = 5500
$ws.Cells.Item(5, 14).Value = "`$/kilo"
$ws.Cells.Item(5, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(5, 16).Value = 5500
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = "Hortaliza"

# Insert a second new data row at row 9 (after the first insert has already
# shifted everything down once), for the new weekly price observation dated 44669.
$ws.Rows("9:9").Insert()
$ws.Cells.Item(9, 1).Value = 10
$ws.Cells.Item(9, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(9, 3).Value = "La Araucanía"
$ws.Cells.Item(9, 4).Value = 44669
$ws.Cells.Item(9, 5).Value = 9
$ws.Cells.Item(9, 6).Value = 100112042
$ws.Cells.Item(9, 7).Value = "Locoto"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 60
$ws.Cells.Item(9, 11).Value = 6250
$ws.Cells.Item(9, 12).Value = 6250
$ws.Cells.Item(9, 13).Value = 6250
$ws.Cells.Item(9, 14).Value = "`$/kilo"
$ws.Cells.Item(9, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(9, 16).Value = 6250
$ws.Cells.Item(9, 17).Value = 1
$ws.Cells.Item(9, 18).Value = "Hortaliza"
